$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 - Tên đăng nhập
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "50"
$ws.Range("G5").Value = "Không chứa ký tự đặc biệt"
$ws.Range("H5").Value = "`"user123`""

# Row 6 - Mật khẩu
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "20"
$ws.Range("F6").NumberFormat = "@"
$ws.Range("F6").Value = "8"
$ws.Range("G6").Value = "Phải có ít nhất 1 ký tự số và 1 ký tự đặc biệt"
$ws.Range("H6").Value = "`"P@ssw0rd!`""

# Row 7 - Nhớ mật khẩu
$ws.Range("G7").Value = "Chọn nếu muốn lưu mật khẩu trên trình duyệt"
$ws.Range("H7").Value = "true/false"

# Row 8 - was "Quên mật khẩu", becomes "Đăng nhập"
$ws.Range("B8").Value = "Đăng nhập"
$ws.Range("C8").Value = "Có"
$ws.Range("D8").Value = "Button"
$ws.Range("G8").Value = "Click để xác thực thông tin đăng nhập"
$ws.Range("H8").Value = "N/A"

# Row 9 - was "Đăng nhập", becomes "Quên mật khẩu"
$ws.Range("B9").Value = "Quên mật khẩu"
$ws.Range("C9").Value = "Không"
$ws.Range("D9").Value = "Link"
$ws.Range("G9").Value = "Chuyển đến trang khôi phục mật khẩu"
$ws.Range("H9").Value = "N/A"

# New row 10 - Đăng ký tài khoản
$ws.Range("A10").NumberFormat = "@"
$ws.Range("A10").Value = "6"
$ws.Range("B10").Value = "Đăng ký tài khoản"
$ws.Range("C10").Value = "Không"
$ws.Range("D10").Value = "Link"
$ws.Range("E10").Value = "N/A"
$ws.Range("F10").Value = "N/A"
$ws.Range("G10").Value = "Chuyển đến trang đăng ký tài khoản"
$ws.Range("H10").Value = "N/A"
